$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 1092, shifting existing rows 1092:1141 down to 1093:1142
$ws.Rows.Item(1092).Insert()

# Populate the newly inserted row 1092 with the new record's data
$ws.Cells.Item(1092, 1).Value = 3
$ws.Cells.Item(1092, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(1092, 3).Value = "Coquimbo"
$ws.Cells.Item(1092, 4).Value = 45147
$ws.Cells.Item(1092, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(1092, 5).Value = 5
$ws.Cells.Item(1092, 6).Value = 100112024
$ws.Cells.Item(1092, 7).Value = "Choclo"
$ws.Cells.Item(1092, 8).Value = "Dulce o Americano"
$ws.Cells.Item(1092, 9).Value = "Primera"
$ws.Cells.Item(1092, 10).Value = 80
$ws.Cells.Item(1092, 11).Value = 40000
$ws.Cells.Item(1092, 12).Value = 41000
$ws.Cells.Item(1092, 13).Value = 40438
$ws.Cells.Item(1092, 14).Value = "$/malla 70 unidades"
$ws.Cells.Item(1092, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(1092, 16).Value = 578
$ws.Cells.Item(1092, 17).Value = 70
$ws.Cells.Item(1092, 18).Value = "Hortaliza"
